# Turn off economic retirements for hydro.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("CRpUNL")

# Add explanatory text rows on the About sheet (rows 13-14)
$wsAbout.Range("A13").Value = "Hydro plants typically have other revenue streams and other considerations when determining whether or not to retire them,"
$wsAbout.Range("A14").Value = "so we assume a value of zero"

# Set hydro retirement fraction (B6) to 0
$wsData.Range("B6").Value = 0

# Apply gray fill (White, Background 1, Darker 25% == RGB D9D9D9) to the cells
# that represent "turned off" / zeroed retirement rows
$grayCells = @("B6", "B16", "B17", "B19", "B20", "B21", "B22", "B23", "B24", "B25")
foreach ($addr in $grayCells) {
    $wsData.Range($addr).Interior.Color = 14277081
}

$wsData.Range("G16").Select()
$wsAbout.Range("A15").Select()
